# Temporary adjustment for RTMF
$wb = $excel.ActiveWorkbook
$wsAbout = $wb.Worksheets.Item("About")
$ws = $wb.Worksheets.Item("RTMF-passengers")

# Row 2 = LDVs (the "shifted-from" mode). Update the HDVs and rail
# fractions, and replace the "Non-motorized/eliminated" formula with its
# resulting static value.
$ws.Range("C2").Value = 0.2
$ws.Range("E2").Value = 0
$ws.Range("I2").Value = 0.8

# Record the last active selection on the RTMF-passengers sheet, then
# restore "About" as the active/selected tab (matching the saved view).
$ws.Range("D7").Select()
$wsAbout.Select()
